# edit.ps1 - apply changes described by the diff to login.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add new cell B2 with value "jxbsansanas" (new shared string)
$ws.Range("B2").Value = "jxbsansanas"

# 2. Set style of B2: center alignment (no special font) -- create this xf first
$ws.Range("B2").HorizontalAlignment = -4108  # xlCenter

# 3. Update style of A2: add horizontal center alignment to its existing vertical-center style
$ws.Range("A2").HorizontalAlignment = -4108  # xlCenter

# 4. Update the selection to C6
$ws.Range("C6").Select()
